$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 98 with latest Argent (Solar) price data.
# Leading apostrophes force these to be stored as text, matching the
# existing rows in the sheet (which store every column, including
# numeric-looking ones, as text).
$ws.Range("A98").Value = "'2025-06-07"
$ws.Range("B98").Value = "'35.5"
$ws.Range("C98").Value = "'35.21"
$ws.Range("D98").Value = "'0.94"
$ws.Range("E98").Value = "'0.248"
$ws.Range("F98").Value = "'0.09"
$ws.Range("G98").Value = "'5,694"
$ws.Range("H98").Value = "'8,526"
$ws.Range("I98").Value = "'8,576"
$ws.Range("J98").Value = "'7.1965"
